$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 900
$ws.Range("B3").Value = 450
$ws.Range("B4").Value = 100
$ws.Range("B5").Value = 200
$ws.Range("B6").Value = 850
